# Apply updated crypto price/volume data (Mon Mar 27 21:42:55 UTC 2023 GitHub Actions run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.149.83"
$ws.Range("E2").Value = "  -2.49%  "
$ws.Range("D3").Value = "1.710.76"
$ws.Range("E3").Value = "  -2.98%  "
$ws.Range("E4").Value = "  -0.19%  "
$ws.Range("D5").Value = "'307.66"
$ws.Range("E5").Value = "  -6.28%  "
$ws.Range("D6").Value = "'0.9999"
$ws.Range("E6").Value = "  -0.14%  "
$ws.Range("D7").Value = "'0.4789"
$ws.Range("E7").Value = "  +7.36%  "
$ws.Range("E8").Value = "  -2.98%  "
$ws.Range("D9").Value = "'41.97"
$ws.Range("E9").Value = "  +0.16%  "
$ws.Range("D10").Value = "'0.07285"
$ws.Range("E10").Value = "  -1.58%  "
$ws.Range("D11").Value = "'1.050"
$ws.Range("E11").Value = "  -4.45%  "
$ws.Range("D12").Value = "'0.9998"
$ws.Range("E12").Value = "  -0.19%  "
$ws.Range("D13").Value = "'19.87"
$ws.Range("E13").Value = "  -4.73%  "
$ws.Range("D14").Value = "'5.853"
$ws.Range("E14").Value = "  -2.71%  "
$ws.Range("D15").Value = "1.708.82"
$ws.Range("E15").Value = "  -3.02%  "
$ws.Range("D16").Value = "'6.843"
$ws.Range("E16").Value = "  -5.44%  "
$ws.Range("D17").Value = "'89.16"
$ws.Range("E17").Value = "  -3.99%  "
$ws.Range("D18").Value = "'0.00001039"
$ws.Range("E18").Value = "  -2.04%  "
$ws.Range("E19").Value = "  -1.26%  "
$ws.Range("D20").Value = "'0.9995"
$ws.Range("D21").Value = "'16.47"
$ws.Range("E21").Value = "  -3.51%  "
$ws.Range("E22").Value = "  -2.88%  "
$ws.Range("D23").Value = "27.186.79"
$ws.Range("E23").Value = "  -2.47%  "
$ws.Range("D24").Value = "'10.88"
$ws.Range("E24").Value = "  -3.16%  "
$ws.Range("D25").Value = "'2.097"
$ws.Range("E25").Value = "  -0.51%  "
$ws.Range("D26").Value = "'154.24"
$ws.Range("E26").Value = "  -3.94%  "
$ws.Range("D27").Value = "'19.68"
$ws.Range("E27").Value = "  -3.23%  "
$ws.Range("D28").Value = "1.902.80"
$ws.Range("E28").Value = "  -3.22%  "
$ws.Range("D29").Value = "'2.080"
$ws.Range("E29").Value = "  -2.88%  "
$ws.Range("D30").Value = "'119.57"
$ws.Range("E30").Value = "  -3.62%  "
$ws.Range("E31").Value = "  -8.08%  "
$ws.Range("D32").Value = "'0.09281"
$ws.Range("E32").Value = "  +0.96%  "
$ws.Range("D33").Value = "'3.580"
$ws.Range("E33").Value = "  -3.03%  "
$ws.Range("D34").Value = "'5.308"
$ws.Range("E34").Value = "  -5.96%  "
$ws.Range("E35").Value = "  -3.59%  "
$ws.Range("D36").Value = "'0.05864"
$ws.Range("E36").Value = "  -5.10%  "
$ws.Range("D37").Value = "'11.06"
$ws.Range("E37").Value = "  -6.46%  "
$ws.Range("D38").Value = "'0.1993"
$ws.Range("E38").Value = "  -4.98%  "
$ws.Range("D39").Value = "'4.748"
$ws.Range("E39").Value = "  -3.97%  "
$ws.Range("B40").Value = "WEMIXTOKEN"
$ws.Range("C40").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D40").Value = "'1.403"
$ws.Range("E40").Value = "  +0.80%  "
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").Value = "'0.5898"
$ws.Range("E41").Value = "  -6.49%  "
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").Value = "'1.109"
$ws.Range("E42").Value = "  -6.20%  "
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").Value = "'7.449"
$ws.Range("E43").Value = "  -4.86%  "
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").Value = "'12.65"
$ws.Range("E44").Value = "  -3.90%  "
$ws.Range("B45").Value = "PancakeSwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D45").Value = "'3.558"
$ws.Range("E45").Value = "  -4.79%  "
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").Value = "'0.5619"
$ws.Range("E46").Value = "  -3.79%  "
$ws.Range("B47").Value = "Quant"
$ws.Range("C47").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D47").Value = "'117.82"
$ws.Range("E47").Value = "  -3.64%  "
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").Value = "'1.842"
$ws.Range("E48").Value = "  -5.51%  "
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value = "'0.06625"
$ws.Range("E49").Value = "  -3.69%  "
$ws.Range("B50").Value = "EOS"
$ws.Range("C50").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D50").Value = "'1.085"
$ws.Range("E50").Value = "  -4.33%  "
$ws.Range("B51").Value = "PaxDollar"
$ws.Range("C51").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D51").Value = "'0.9997"
$ws.Range("E51").Value = "  -0.11%  "
